$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells that would otherwise be parsed as numbers
$priceTextCells = @("D4", "D5", "D6", "D7", "D10", "D12", "D13", "D14", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D41", "D42", "D44", "D46", "D47", "D49", "D51")
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '66.921.19'
$ws.Range('E2').Value = '  +6.80%  '
$ws.Range('D3').Value = '3.549.32'
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '413.91'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').Value = '128.83'
$ws.Range('E6').Value = '  -1.16%  '
$ws.Range('D7').Value = '0.646'
$ws.Range('E7').Value = '  +3.74%  '
$ws.Range('D8').Value = '3.542.31'
$ws.Range('E8').Value = '  +2.34%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '0.774'
$ws.Range('E10').Value = '  +6.55%  '
$ws.Range('E11').Value = '  +25.71%  '
$ws.Range('D12').Value = '0.0000330'
$ws.Range('E12').Value = '  +50.44%  '
$ws.Range('D13').Value = '42.28'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = '9.86'
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').Value = '4.124.08'
$ws.Range('E15').Value = '  +2.68%  '
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '3.636.36'
$ws.Range('E17').Value = '  +4.77%  '
$ws.Range('D18').Value = '20.12'
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('E19').Value = '  +2.94%  '
$ws.Range('D20').Value = '66.858.32'
$ws.Range('E20').Value = '  +6.74%  '
$ws.Range('D21').Value = '12.32'
$ws.Range('E21').Value = '  -3.55%  '
$ws.Range('D22').Value = '444.94'
$ws.Range('E22').Value = '  -5.42%  '
$ws.Range('D23').Value = '89.15'
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('D24').Value = '3.12'
$ws.Range('E24').Value = '  -4.42%  '
$ws.Range('D25').Value = '12.91'
$ws.Range('E25').Value = '  -3.75%  '
$ws.Range('D26').Value = '3.31'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').Value = '9.95'
$ws.Range('E27').Value = '  -6.00%  '
$ws.Range('D28').Value = '34.46'
$ws.Range('E28').Value = '  +3.39%  '
$ws.Range('D29').Value = '4.85'
$ws.Range('E29').Value = '  +0.98%  '
$ws.Range('E30').Value = '  +3.71%  '
$ws.Range('D31').Value = '12.25'
$ws.Range('E31').Value = '  +2.06%  '
$ws.Range('E32').Value = '  +2.89%  '
$ws.Range('D33').Value = '7.25'
$ws.Range('E33').Value = '  -4.70%  '
$ws.Range('D34').Value = '0.157'
$ws.Range('E34').Value = '  -6.03%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = '39.08'
$ws.Range('E36').Value = '  -4.96%  '
$ws.Range('D37').Value = '56.26'
$ws.Range('E37').Value = '  -4.23%  '
$ws.Range('D38').Value = '0.0489'
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').Value = '0.0₃0730'
$ws.Range('E39').Value = '  +28.01%  '
$ws.Range('E40').Value = '  +9.30%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').Value = '146.84'
$ws.Range('E42').Value = '  +0.97%  '
$ws.Range('E43').Value = '  -3.94%  '
$ws.Range('D44').Value = '2.69'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('E45').Value = '  -2.04%  '
$ws.Range('D46').Value = '3.19'
$ws.Range('E46').Value = '  -5.01%  '
$ws.Range('D47').Value = '0.304'
$ws.Range('E47').Value = '  -5.43%  '
$ws.Range('E48').Value = '  -6.35%  '
$ws.Range('D49').Value = '117.38'
$ws.Range('E49').Value = '  +7.18%  '
$ws.Range('E50').Value = '  -6.61%  '
$ws.Range('B51').Value = 'Celestia'
$ws.Range('C51').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D51').Value = '15.31'
$ws.Range('E51').Value = '  -6.53%  '
